$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - new columns I (I0) and J (IF), copying the header style from H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I (I0) and J (IF) per row
$data = @{
    2  = @(1, 5)
    3  = @(1, 7)
    4  = @(1, 7)
    5  = @(1, 4)
    6  = @(1, 4)
    7  = @(1, 6)
    8  = @(1, 5)
    9  = @(1, 4)
    10 = @(1, 4)
    11 = @(1, 3)
    12 = @(1, 4)
    13 = @(1, 5)
    14 = @(1, 6)
    15 = @(1, 6)
    16 = @(1, 5)
    17 = @(1, 6)
    18 = @(1, 5)
    19 = @(7, 8)
    20 = @(7, 9)
    21 = @(1, 3)
    22 = @(1, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
